# Updates the cryptos price/volume table (columns D and E) with the
# latest scraped values. D/E cells are stored as plain text in the sheet
# (prices use "." as a thousands-style separator, not a decimal point, so
# they must stay literal text). Values that otherwise look like valid
# numbers (e.g. "1.001") are written with a leading apostrophe so Excel
# keeps them as text instead of silently coercing them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.872.41"
$ws.Range("E2").Value = "  +4.21%  "
$ws.Range("D3").Value = "1.879.98"
$ws.Range("E3").Value = "  +3.55%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'278.43"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5315"
$ws.Range("E7").Value = "  +4.09%  "
$ws.Range("D8").Value = "'0.3454"
$ws.Range("E8").Value = "  -1.80%  "
$ws.Range("D9").Value = "'45.06"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").Value = "'0.06969"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").Value = "'20.11"
$ws.Range("E11").Value = "  +0.45%  "
$ws.Range("D12").Value = "'0.8061"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").Value = "'0.07734"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "1.882.85"
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("D15").Value = "'90.54"
$ws.Range("E15").Value = "  +3.56%  "
$ws.Range("D16").Value = "'5.181"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "'14.57"
$ws.Range("E17").Value = "  +3.17%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "'0.000008038"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "26.931.15"
$ws.Range("E21").Value = "  +4.15%  "
$ws.Range("D22").Value = "2.123.05"
$ws.Range("E22").Value = "  +3.48%  "
$ws.Range("D23").Value = "'4.756"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").Value = "'10.05"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("D25").Value = "'6.215"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "'2.367"
$ws.Range("E26").Value = "  +8.01%  "
$ws.Range("D27").Value = "'147.07"
$ws.Range("E27").Value = "  +4.19%  "
$ws.Range("D28").Value = "'1.662"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "'17.36"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").Value = "'113.60"
$ws.Range("E30").Value = "  +3.75%  "
$ws.Range("D31").Value = "'4.360"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'4.324"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("D33").Value = "'0.08895"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("D34").Value = "'0.04937"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "'1.175"
$ws.Range("E35").Value = "  +3.47%  "
$ws.Range("D36").Value = "'0.7333"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'2.898"
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").Value = "'3.295"
$ws.Range("E38").Value = "  +4.57%  "
$ws.Range("D39").Value = "'2.375"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "'0.01854"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "'0.5146"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").Value = "'0.9576"
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").Value = "'115.71"
$ws.Range("E43").Value = "  +3.60%  "
$ws.Range("D44").Value = "'6.191"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").Value = "'8.133"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.04%  "
$ws.Range("D47").Value = "'0.4478"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").Value = "'0.1344"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").Value = "'9.396"
$ws.Range("E49").Value = "  +1.88%  "
$ws.Range("D50").Value = "'36.24"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("D51").Value = "'0.05951"
$ws.Range("E51").Value = "  +1.91%  "
